$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.467.35"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +2.11%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.729.48"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +2.76%  "

$ws.Range("E4").Value = "  +0.27%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.96"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.54%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.0000"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.20%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4802"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +3.34%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2678"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +2.56%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06234"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.23%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.730.37"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +3.29%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07134"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.74%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.75"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +4.99%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6181"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +6.27%  "

$ws.Range("E14").Value = "  +4.06%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.20"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.08%  "

$ws.Range("E16").Value = "  +0.17%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.478.05"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.21%  "

$ws.Range("E18").Value = "  +0.16%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006961"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +3.60%  "

$ws.Range("E20").Value = "  +2.44%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.954.73"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.46%  "

$ws.Range("E22").Value = "  +1.73%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.935"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.80%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.324"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.58%  "

$ws.Range("E25").Value = "  +1.77%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.37"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.30%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.797"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +4.49%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.404"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.88%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "106.76"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.78%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.987"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.72%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08023"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +3.95%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.745"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +3.13%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04561"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +4.25%  "

$ws.Range("E34").Value = "  +0.83%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6422"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +5.54%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9929"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.65%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9445"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.42%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.995"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +6.04%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "108.05"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.20%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.408"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.72%  "

$ws.Range("E41").Value = "  +0.91%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01504"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +3.04%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.678"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +11.96%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3923"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +4.97%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.968"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +12.63%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1195"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +6.44%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05325"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.45%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "30.86"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.15%  "

$ws.Range("E49").Value = "  +2.43%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.274"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +4.90%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3438"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +3.28%  "
